$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.107.63"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.960.83"
$ws.Range("E3").Value = "  +0.25%  "
$__style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $__style
$ws.Range("E4").Value = "  +0.01%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.16"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +1.43%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.29"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("E8").Value = "  +0.00%  "
$__style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = $__style
$ws.Range("E9").Value = "  +1.08%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.54"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("E11").Value = "  -1.33%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$__style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.80"
$ws.Range("D13").Style = $__style
$ws.Range("E13").Value = "  +5.41%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.418.74"
$ws.Range("E14").Value = "  +0.06%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.36"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  +2.21%  "
$__style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.50"
$ws.Range("D16").Style = $__style
$ws.Range("E16").Value = "  +27.78%  "
$ws.Range("D17").Value = "2.952.70"
$ws.Range("E17").Value = "  +0.60%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.999"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "51.180.72"
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.13"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  -0.90%  "
$__style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = $__style
$ws.Range("E21").Value = "  -2.23%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.40%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.33"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +13.83%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.18"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +2.54%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "266.95"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +0.80%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.88"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  -7.02%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.22"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  -10.12%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  -2.74%  "
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.87"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  +0.72%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.111"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -2.38%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.31"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +4.14%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.36"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.12"
$ws.Range("D34").Style = $__style
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +2.00%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0436"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -0.09%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  +7.95%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.117"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$__style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "16.54"
$ws.Range("D41").Style = $__style
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  -2.94%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.79"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.53"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +6.98%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.44"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.272"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  -5.94%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$__style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("D47").Style = $__style
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("D49").Value = "2.050.93"
$ws.Range("E49").Value = "  +3.44%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0321"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  -7.68%  "
$ws.Range("E51").Value = "  +6.56%  "
